$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Page 3" of results: rows 26-32 shift up by one (each row now shows what
# used to be the next row's data), and row 33 becomes a brand-new entry.
$data = @(
    @{ Row = 26; A = "BSCI ISO9001 tessuto di bambù salviette di bambù lavabili all'ingrosso 100% rotolo di asciugamani di carta di bambù riutilizzabile da cucina"; B = "0,2166 €"; C = "Ordine minimo: 10.000 sacchi"; D = "Ningbo Riway Nonwovens Tech Co., Ltd."; E = "4.7" },
    @{ Row = 27; A = "Asciugamani da cucina riutilizzabili ecologici personalizzati con Design goffrato e assorbente per pulire il rotolo da cucina"; B = "0,5197 €"; C = "Ordine minimo: 1 rullo"; D = "Hangzhou Shengbo Cleaning Product Co., Ltd."; E = "5.0" },
    @{ Row = 28; A = "Asciugamani di carta velina pulita da cucina piegati riutilizzabili riciclati rotolo di asciugamani di carta da cucina ad assorbimento di olio forte usa e getta"; B = "0,3465-0,4504 €"; C = "Ordine minimo: 1.000 rulli"; D = "Hangzhou Micker Sanitary Products Co., Ltd."; E = "4.8" },
    @{ Row = 29; A = "Asciugamano di carta riutilizzabile,"; B = "8,58 €"; C = "Ordine minimo: 10 parti"; D = "Zhejiang Chloven Cosmetics Co., Ltd."; E = $null },
    @{ Row = 30; A = "Carta usa e getta eco-friendly grande rotolo forte pulizia riutilizzabile straccio pigro campione di carta da cucina"; B = "0,0867-0,2599 €"; C = "Ordine minimo: 20.000 parti"; D = "Hangzhou Lin'an Sanxin Cleaning Products Co., Ltd."; E = "__SKIP__" },
    @{ Row = 31; A = "Pasta di legno vergine di alta qualità usa e getta carta da cucina asciugamani morbidi e assorbenti lavabili Private Label per uso alberghiero"; B = "0,13-0,2079 €"; C = "Ordine minimo: 12.000 sacchi"; D = "Hebei Yihoucheng Commodity Co., Ltd."; E = "5.0" },
    @{ Row = 32; A = "Pulizia della cucina Del Prodotto Senza Carta Tovagliolo Di Bambù Organico Unpaper Asciugamani Riutilizzabili"; B = "0,3292 €"; C = "Ordine minimo: 500 parti"; D = "Yiwu Yozeal Commodity Co., Ltd."; E = "4.8" },
    @{ Row = 33; A = "Logo personalizzato cucina asciugamani di carta riutilizzabili lavabili rotolo spesso due strati Ultra assorbente asciugamani in cotone organico senza carta"; B = "0,3898 €"; C = "Ordine minimo: 500 parti"; D = "Yiwu Xinyao Trading Co., Ltd."; E = "5.0" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D

    if ($item.E -eq "__SKIP__") {
        # Already blank in both the old and new data for this row - leave it
        # untouched so the existing (blank) cell isn't disturbed.
        continue
    }

    $eCell = $ws.Range("E$r")
    if ($item.E) {
        # Ratings like "4.7" / "5.0" must stay text (matching the rest of the
        # column), so force text formatting before writing the value -
        # otherwise Excel would coerce it into the number 4.7 / 5.
        $eCell.NumberFormat = "@"
        $eCell.Value = $item.E
    } else {
        $eCell.Value = ""
    }
}
